# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list on Fri Jun 23 06:18:58 UTC 2023 with GitHub Actions".
# Columns: A=index (unchanged), B=Coin, C=Link, D=Price, E=Volume(1h).
# A handful of rows (25,26,34,35,48,49,50) also had their Coin/Link swapped
# with the neighboring row as the ranking shifted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $text) {
    # Force the cell to stay a text value (many of these "prices" look like
    # numbers, e.g. 0.9982, but must remain text as in the source data,
    # matching values such as 30.007.98 that are not valid numbers anyway).
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "30.007.98"
$ws.Range("E2").Value = "  -0.48%  "

Set-TextValue $ws "D3" "1.882.60"
$ws.Range("E3").Value = "  -1.36%  "

Set-TextValue $ws "D4" "0.9982"
$ws.Range("E4").Value = "  -0.41%  "

Set-TextValue $ws "D5" "243.18"
$ws.Range("E5").Value = "  -3.09%  "

Set-TextValue $ws "D6" "0.9983"
$ws.Range("E6").Value = "  -0.35%  "

Set-TextValue $ws "D7" "0.4966"
$ws.Range("E7").Value = "  -2.60%  "

Set-TextValue $ws "D8" "0.2928"
$ws.Range("E8").Value = "  -0.53%  "

Set-TextValue $ws "D9" "0.06654"
$ws.Range("E9").Value = "  -1.44%  "

Set-TextValue $ws "D10" "1.880.83"
$ws.Range("E10").Value = "  -1.45%  "

Set-TextValue $ws "D11" "16.79"
$ws.Range("E11").Value = "  -2.54%  "

Set-TextValue $ws "D12" "0.07246"
$ws.Range("E12").Value = "  -1.21%  "

Set-TextValue $ws "D13" "0.6686"
$ws.Range("E13").Value = "  -2.79%  "

Set-TextValue $ws "D14" "86.41"
$ws.Range("E14").Value = "  -0.12%  "

Set-TextValue $ws "D15" "4.918"
$ws.Range("E15").Value = "  +1.56%  "

Set-TextValue $ws "D16" "29.988.10"
$ws.Range("E16").Value = "  -0.60%  "

Set-TextValue $ws "D17" "0.000007920"
$ws.Range("E17").Value = "  -1.62%  "

Set-TextValue $ws "D18" "0.9983"
$ws.Range("E18").Value = "  -0.39%  "

Set-TextValue $ws "D19" "12.80"
$ws.Range("E19").Value = "  -0.92%  "

Set-TextValue $ws "D20" "2.124.76"
$ws.Range("E20").Value = "  -1.62%  "

Set-TextValue $ws "D21" "0.9973"
$ws.Range("E21").Value = "  -0.48%  "

Set-TextValue $ws "D22" "4.787"
$ws.Range("E22").Value = "  -0.51%  "

Set-TextValue $ws "D23" "5.765"
$ws.Range("E23").Value = "  +1.17%  "

Set-TextValue $ws "D24" "9.090"
$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("B25").Value = "BitcoinCash"
$ws.Range("C25").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws "D25" "143.14"
$ws.Range("E25").Value = "  +6.50%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws "D26" "149.90"
$ws.Range("E26").Value = "  +1.89%  "

Set-TextValue $ws "D27" "17.11"
$ws.Range("E27").Value = "  +0.40%  "

Set-TextValue $ws "D28" "1.922"
$ws.Range("E28").Value = "  -3.46%  "

Set-TextValue $ws "D29" "1.390"
$ws.Range("E29").Value = "  -0.28%  "

Set-TextValue $ws "D30" "4.203"
$ws.Range("E30").Value = "  -0.15%  "

Set-TextValue $ws "D31" "0.08783"
$ws.Range("E31").Value = "  +0.24%  "

Set-TextValue $ws "D32" "3.967"
$ws.Range("E32").Value = "  -0.21%  "

Set-TextValue $ws "D33" "0.05085"
$ws.Range("E33").Value = "  +0.56%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws "D34" "1.124"
$ws.Range("E34").Value = "  -1.83%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws "D35" "0.7149"
$ws.Range("E35").Value = "  +0.71%  "

Set-TextValue $ws "D36" "2.665"
$ws.Range("E36").Value = "  -1.07%  "

Set-TextValue $ws "D37" "0.01814"
$ws.Range("E37").Value = "  +8.09%  "

Set-TextValue $ws "D38" "2.694"
$ws.Range("E38").Value = "  -4.27%  "

Set-TextValue $ws "D39" "2.181"
$ws.Range("E39").Value = "  -4.31%  "

Set-TextValue $ws "D40" "0.9331"
$ws.Range("E40").Value = "  -3.56%  "

Set-TextValue $ws "D41" "5.811"
$ws.Range("E41").Value = "  -3.90%  "

Set-TextValue $ws "D42" "0.4265"
$ws.Range("E42").Value = "  +0.06%  "

Set-TextValue $ws "D43" "0.9983"
$ws.Range("E43").Value = "  -0.06%  "

Set-TextValue $ws "D44" "102.33"
$ws.Range("E44").Value = "  -2.55%  "

Set-TextValue $ws "D45" "7.486"
$ws.Range("E45").Value = "  -1.14%  "

Set-TextValue $ws "D46" "0.1268"
$ws.Range("E46").Value = "  -0.14%  "

Set-TextValue $ws "D47" "0.05657"
$ws.Range("E47").Value = "  -1.44%  "

$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws "D48" "0.3819"
$ws.Range("E48").Value = "  +0.92%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws "D49" "8.350"
$ws.Range("E49").Value = "  -0.90%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue $ws "D50" "32.58"
$ws.Range("E50").Value = "  -1.19%  "

Set-TextValue $ws "D51" "56.08"
$ws.Range("E51").Value = "  -0.80%  "
